# Re-run of the NATMI script with updated TPM data for Clec4g-Lag3.
# The sending clusters for this ligand-receptor pair shifted
# (rows 2-5: ECs -> FAPs, rows 6-9: FAPs -> MuSCs), while the target
# clusters (column D) and ligand/receptor symbols (B/C) stay the same.
# The expression / specificity values in columns G,H,I,J,M,N,O,P,Q,R,S,T
# also change due to the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ A="FAPs"; G=0.1416773333333333; H=0.425032; I=0.4649030121259454; J=0.4649030121259454; M=9.278280333333333; N=27.834841; O=0.2122966588143784; P=0.2122966588143784; Q=1.314522015545778; R=11.830698139912; S=0.09869735614707865; T=0.09869735614707865 }
    3 = @{ A="FAPs"; G=0.1416773333333333; H=0.425032; I=0.4649030121259454; J=0.4649030121259454; O=0.2154323368929792; P=0.2154323368929792; Q=1.333937855112; R=12.005440696008; S=0.1001551423308775; T=0.1001551423308774 }
    4 = @{ A="FAPs"; G=0.1416773333333333; H=0.425032; I=0.4649030121259454; J=0.4649030121259454; M=7.033255; N=21.099765; O=0.1609281551588013; P=0.1609281551588013; Q=0.9964528130533336; R=8.96807531748; S=0.07481598406919822; T=0.0748159840691982 }
    5 = @{ A="FAPs"; G=0.1416773333333333; H=0.425032; I=0.4649030121259454; J=0.4649030121259454; M=17.977458; N=53.932374; O=0.4113428491338411; P=0.411342849133841; Q=2.546998309552; R=22.922984785968; S=0.191234529578791; T=0.191234529578791 }
    6 = @{ A="MuSCs"; G=0.1630686666666667; H=0.489206; I=0.5350969878740547; J=0.5350969878740547; M=9.278280333333333; N=27.834841; O=0.2122966588143784; P=0.2122966588143784; Q=1.512996802916222; R=13.616971226246; S=0.1135993026672998; T=0.1135993026672998 }
    7 = @{ A="MuSCs"; G=0.1630686666666667; H=0.489206; I=0.5350969878740547; J=0.5350969878740547; O=0.2154323368929792; P=0.2154323368929792; Q=1.535344167846; R=13.818097510614; S=0.1152771945621018; T=0.1152771945621017 }
    8 = @{ A="MuSCs"; G=0.1630686666666667; H=0.489206; I=0.5350969878740547; J=0.5350969878740547; M=7.033255; N=21.099765; O=0.1609281551588013; P=0.1609281551588013; Q=1.146903515176667; R=10.32213163659; S=0.08611217108960309; T=0.08611217108960308 }
    9 = @{ A="MuSCs"; G=0.1630686666666667; H=0.489206; I=0.5350969878740547; J=0.5350969878740547; M=17.977458; N=53.932374; O=0.4113428491338411; P=0.4113428491338411; Q=2.931560106116; R=26.384040955044; S=0.2201083195550501; T=0.22010831955505 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
